$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 'Подавите любое сопротивление {TEAM_TAR.FactionDef.Name} на территории базы, командир. Оливейра конец связи.'
$ws.Range("C16").Value = 'Уничтожьте оборонительные силы {TEAM_TAR.FactionDef.Name}'
$ws.Range("C19").Value = 'Удерживайте базу до прибытия сил {TEAM_EMP.FactionDef.Name}'
$ws.Range("C22").Value = 'Уничтожьте подкрепления {TEAM_TAR.FactionDef.Name}'
$ws.Range("C24").Value = 'Торговые переговоры между {TEAM_EMP.FactionDef.Name} и {TEAM_TAR.FactionDef.Name} зашли в тупик. Мы считаем, что если {TEAM_TAR.FactionDef.Name} потеряют свою военную базу в системе {TGT_SYSTEM.Name}, это вынудит их снова сесть за стол переговоров. Ваша задача - захватить эту базу и удерживать до подхода наших войск.'
$ws.Range("C26").Value = 'Мы здесь, чтобы захватить исследовательскую станцию ​​{TEAM_TAR.FactionDef.Name}.'
$ws.Range("C43").Value = 'Уничтожьте силы обороны {TEAM_TAR.FactionDef.Name}'
$ws.Range("C49").Value = 'Уничтожьте подкрепления {TEAM_TAR.FactionDef.Name}'
$ws.Range("C51").Value = 'Мы давно подозревали, что в системе {TGT_SYSTEM.name} находится секретная исследовательская станция {TEAM_TAR.FactionDef.Name}, и наконец узнали о ее возможном местонахождении. Нам нужно несколько смелых и молчаливых наемников, чтобы обезопасить это место для прилёта дропшипа и препроводить в него весь ценный научный персонал.'
$ws.Range("C69").Value = 'Уничтожьте гарнизон {TEAM_TAR.FactionDef.Name}'
$ws.Range("C75").Value = 'Уничтожьте подкрепления {TEAM_TAR.FactionDef.Name}'
$ws.Range("C77").Value = 'Мы обнаружили, что {TEAM_TAR.FactionDef.Name} повторно активировали старый производственный объект в этой системе. Чтобы сохранить баланс сил, {TEAM_EMP.FactionDef.Name} должны захватить этот объект. От вас требуется нанести хирургически точный удар, захватить базу и удерживать ее до прихода сил {TEAM_EMP.FactionDef.Name}.'
$ws.Range("C92").Value = 'Ох не люблю я, когда в целях миссии числится «какой-то груз», командир. Но, чем бы он ни был, - пусть лучше он будет у нас, чем в руках {TEAM_TAR.FactionDef.Name}.'
$ws.Range("C107").Value = 'Чудный денёк, командор. А знаете почему? Я только что получил свежие данные о перевозках от наших друзей на {TGT_SYSTEM.name}. Оказывается все ценные грузы {TEAM_TAR.FactionDef.Name}, что планируется отправить с планеты, хранятся на одной конкретной базе. Ну разве не чудо? Там будет очень ценный груз. И мы хотим всё это «одолжить». Вы берёте базу, мы ее грабим, все будут в плюсе.'
$ws.Range("C109").Value = 'Вот координаты, которые мы получили от разведки {TEAM_EMP.FactionDef.Name}.'
$ws.Range("C134").Value = 'Уничтожьте подкрепления {TEAM_TAR.FactionDef.Name}'
$ws.Range("C137").Value = 'Ничем не примечательные пустые системы, подобные этой, - идеальное место для тайных операций по сбору разведданных. Мы полагаем, что здесь, на {TGT_SYSTEM.Name}, есть оборудование {TEAM_TAR.FactionDef.Name}, предназначенное именно для этого. Мониторинг передвижений грузовых кораблей привёл аналитиков {TEAM_EMP.FactionDef.Name} к отдалённому объекту, который все считали заброшенным. Ровно до тех пор, пока всего несколько дней назад мы не обнаружили явную сигнатуру включения термоядерного двигателя меха. Этот объект не заброшен и мы хотим получить данные, которые там хранятся.'

$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("D57").Select()
